# MENT-169: add UATS cabinet
# Appends a new "UATS" entry to the list of cabinets in column A,
# right after the existing "Laboratório" row, and moves the active
# selection the way the authored workbook recorded it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing data rows run from A2:A14 (A1 is the "Name" header).
# The new cabinet goes in the next empty row, A15.
$ws.Range("A15").Value = "UATS"

# Match the formatting of the other cabinet rows (style index 1 in the
# saved file) by copying the format from the row right above it instead
# of setting font properties individually (which this host only commits
# the first of, per range, when chained).
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)

# Reproduce the recorded selection state after the edit.
[void]$ws.Range("A25").Select()
